$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert two new rows (new "line7" / "line8" entries) right after the
# existing "line6" row (row 7), pushing the "extr1..extr8" rows down
# from 8-15 to 10-17.
# ------------------------------------------------------------------
$ws.Range("A8:A9").EntireRow.Insert()

# Row-insert clones the format of the row above but drops the border,
# which would otherwise allocate a stray style for column A. Restore
# the original bordered style used by every other row in column A.
$ws.Range("A10").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Fill in the full data block for rows 8-17 (name, from_bus, to_bus,
# in_service) to match the updated contingency table.
# ------------------------------------------------------------------
$names      = @("line7", "line8", "extr1", "extr2", "extr3", "extr4", "extr5", "extr6", "extr7", "extr8")
$fromBus    = @(14, 16, 5, 5, 10, 7, 9, 7, 5, 8)
$toBus      = @(11, 9, 12, 9, 11, 8, 11, 11, 7, 5)
$inService  = @($true, $true, $true, $true, $false, $true, $true, $true, $true, $true)

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = 8 + $i
    $ws.Cells.Item($r, 1).Value = 6 + $i
    $ws.Cells.Item($r, 2).Value = $names[$i]
    $ws.Cells.Item($r, 3).Value = $fromBus[$i]
    $ws.Cells.Item($r, 4).Value = $toBus[$i]
    $ws.Cells.Item($r, 5).Value = $inService[$i]
}
